# Update countries & provincias Spain
#
# Refreshes the "Pais" COVID dashboard sheet: bumps the "last updated"
# timestamp banner and pushes new per-country figures (Casos totales,
# Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes)
# for the rows whose source numbers moved since the previous snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Banner in A1: "Datos actualizados a ..."
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 13:11"

# Row -> country -> { column letter : new value }
# Columns: B Casos totales, C Nuevos casos, D Casos activos,
#          E Recuperados, F Casos criticos, G Muertes hoy, H Muertes
$updates = @{
    15  = @{ B = 393425; C = 2313; D = 339111; E = 31645;          G = 127; H = 22669 }  # Iran
    37  = @{ B = 98304;  C = 1271; D = 41002;  E = 53284;          G = 51;  H = 4018  }  # Rumania
    45  = @{ B = 75981;  C = 883;  D = 67359;  E = 8229;           G = 2;   H = 393   }  # Marruecos
    46  = @{ B = 75721;            D = 57239;  E = 17055;                   H = 1427  }  # Emiratos Arabes Unidos
    58  = @{ B = 49219;  C = 1081; D = 33882;  E = 15025;          G = 6;   H = 312   }  # Nepal
    60  = @{ B = 45306;  C = 469;  D = 37700;  E = 5588;           G = 0;   H = 2018  }  # Armenia
    61  = @{ B = 45152;  C = 199;  D = 41023;  E = 3224;           G = 2;   H = 905   }  # Ghana
    62  = @{ B = 45012;            D = 43898;  E = 831;                    H = 283   }  # Suiza
    87  = @{ B = 14102;  C = 58;   D = 10176;  E = 3633;           G = 1;   H = 293   }  # Senegal
    145 = @{ B = 2162;   C = 63;   D = 1760;   E = 388;                    H = 14    }  # Islandia
    146 = @{ B = 2150;             D = 2060;   E = 80;                     H = 10    }  # Botsuana
    147 = @{ B = 2126;             D = 493;    E = 1624;                   H = 9     }  # Malta
    165 = @{ B = 1059;  C = 5;     D = 890;    E = 134                              }  # Vietnam
    183 = @{ B = 322;   C = 2;     D = 277;    E = 45                               }  # Gibraltar
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
